$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pas_005")
$ws.Columns("B:F").Hidden = $false
$ws.Columns("B:F").EntireColumn.AutoFit()
Write-Host "done"
